$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text storage (values look numeric but are stored as text in the source)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.372.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6265"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07434"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2897"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.86"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07715"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.837.39"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.957"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6747"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001021"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.72"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.217"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.348.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "232.02"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.373"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.467"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1345"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07262"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +12.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.461"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.478"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.043"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.050"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.817"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.139"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6940"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.571"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01837"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.907"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.233.05"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9474"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9991"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.985.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.60"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.48"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.716"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000116"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.950"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.909"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1138"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3902"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.22%  "
